# "Se agrega el campo para EGMA"
# Add a new "egma_items" column as the last column on every sheet of the
# gradebook workbook, then leave the workbook/sheets with the selection
# state the author ended up with after making the edit (4th tab, "5°",
# active).

$wb = $excel.ActiveWorkbook

$sheetNames = @("2°", "3°", "4°", "5°")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Find the last used column in row 1 (the header row) and append the
    # new header right after it.
    $lastCol = $ws.Cells(1, $ws.Columns.Count).End(-4159).Column
    $newCol = $lastCol + 1

    $ws.Cells(1, $newCol).Value = "egma_items"
}

# Restore each sheet's own selection.
$null = $wb.Worksheets.Item("2°").Range("X8:X9").Select()
$null = $wb.Worksheets.Item("3°").Range("Y1").Select()
$null = $wb.Worksheets.Item("4°").Range("W1").Select()
$null = $wb.Worksheets.Item("5°").Range("W10").Select()

# "5°" ends up being the active tab.
$null = $wb.Worksheets.Item("5°").Activate()
